$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 466.66666  # H2
$ws.Cells.Item(2, 9).Value = 500  # I2
$ws.Cells.Item(2, 10).Value = 450  # J2
$ws.Cells.Item(2, 11).Value = 500  # K2
$ws.Cells.Item(2, 12).Value = 450  # L2
$ws.Cells.Item(2, 13).Value = -387  # M2
$ws.Cells.Item(2, 14).Value = -676  # N2
$ws.Cells.Item(15, 8).Value = 226404.05  # H15
$ws.Cells.Item(15, 9).Value = 226404.05  # I15
$ws.Cells.Item(15, 11).Value = 679212.1499999999  # K15
$ws.Cells.Item(15, 13).Value = -679043.1499999999  # M15
$ws.Cells.Item(100, 8).Value = 3478.2144  # H100
$ws.Cells.Item(100, 9).Value = 2355.111  # I100
$ws.Cells.Item(100, 10).Value = 5499.8  # J100
$ws.Cells.Item(100, 11).Value = 2355.111  # K100
$ws.Cells.Item(100, 12).Value = 5499.8  # L100
$ws.Cells.Item(100, 13).Value = -1814.111  # M100
$ws.Cells.Item(100, 14).Value = -6581.8  # N100
$ws.Cells.Item(131, 8).Value = 2639.7  # H131
$ws.Cells.Item(131, 9).Value = 924.625  # I131
$ws.Cells.Item(131, 10).Value = 9500  # J131
$ws.Cells.Item(131, 11).Value = 2773.875  # K131
$ws.Cells.Item(131, 12).Value = 28500  # L131
$ws.Cells.Item(131, 13).Value = 2266.125  # M131
$ws.Cells.Item(131, 14).Value = -38580  # N131
$ws.Cells.Item(138, 8).Value = 5623.7573  # H138
$ws.Cells.Item(138, 10).Value = 5645.5713  # J138
$ws.Cells.Item(138, 12).Value = 16936.7139  # L138
$ws.Cells.Item(138, 14).Value = -27216.7139  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1810.8889  # H5
$ws.Cells.Item(5, 9).Value = 1833.1666  # I5
$ws.Cells.Item(5, 10).Value = 1766.3334  # J5
$ws.Cells.Item(5, 11).Value = 1833.1666  # K5
$ws.Cells.Item(5, 12).Value = 1766.3334  # L5
$ws.Cells.Item(5, 13).Value = -1721.1666  # M5
$ws.Cells.Item(5, 14).Value = -1990.3334  # N5
$ws.Cells.Item(32, 8).Value = 2610150.2  # H32
$ws.Cells.Item(32, 9).Value = 3289564  # I32
$ws.Cells.Item(32, 10).Value = 19884.812  # J32
$ws.Cells.Item(32, 11).Value = 3289564  # K32
$ws.Cells.Item(32, 12).Value = 19884.812  # L32
$ws.Cells.Item(32, 13).Value = -3289277  # M32
$ws.Cells.Item(32, 14).Value = -20458.812  # N32
$ws.Cells.Item(61, 8).Value = 3100.8857  # H61
$ws.Cells.Item(61, 9).Value = 2621.724  # I61
$ws.Cells.Item(61, 11).Value = 2621.724  # K61
$ws.Cells.Item(61, 13).Value = -2409.724  # M61
$ws.Cells.Item(63, 8).Value = 3231.0625  # H63
$ws.Cells.Item(63, 10).Value = 3679.8  # J63
$ws.Cells.Item(63, 12).Value = 3679.8  # L63
$ws.Cells.Item(63, 14).Value = -5051.8  # N63
$ws.Cells.Item(66, 8).Value = 3231.0625  # H66
$ws.Cells.Item(66, 10).Value = 3679.8  # J66
$ws.Cells.Item(66, 12).Value = 18399  # L66
$ws.Cells.Item(66, 14).Value = -25263  # N66
$ws.Cells.Item(74, 8).Value = 3979.6667  # H74
$ws.Cells.Item(74, 9).Value = 2107.5334  # I74
$ws.Cells.Item(74, 11).Value = 2107.5334  # K74
$ws.Cells.Item(74, 13).Value = -1233.5334  # M74
$ws.Cells.Item(77, 8).Value = 3979.6667  # H77
$ws.Cells.Item(77, 9).Value = 2107.5334  # I77
$ws.Cells.Item(77, 11).Value = 10537.667  # K77
$ws.Cells.Item(77, 13).Value = -6169.666999999999  # M77
$ws.Cells.Item(104, 8).Value = 60611.668  # H104
$ws.Cells.Item(104, 10).Value = 60611.668  # J104
$ws.Cells.Item(104, 12).Value = 60611.668  # L104
$ws.Cells.Item(104, 14).Value = -67599.66800000001  # N104
$ws.Cells.Item(132, 8).Value = 671219.3  # H132
$ws.Cells.Item(132, 10).Value = 4999.3335  # J132
$ws.Cells.Item(132, 12).Value = 14998.0005  # L132
$ws.Cells.Item(132, 14).Value = -20058.0005  # N132
$ws.Cells.Item(136, 8).Value = 3100.8857  # H136
$ws.Cells.Item(136, 9).Value = 2621.724  # I136
$ws.Cells.Item(136, 11).Value = 7865.172  # K136
$ws.Cells.Item(136, 13).Value = -5315.172  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1810.8889  # H4
$ws.Cells.Item(4, 9).Value = 1833.1666  # I4
$ws.Cells.Item(4, 10).Value = 1766.3334  # J4
$ws.Cells.Item(4, 11).Value = 1833.1666  # K4
$ws.Cells.Item(4, 12).Value = 1766.3334  # L4
$ws.Cells.Item(4, 13).Value = -1718.1666  # M4
$ws.Cells.Item(4, 14).Value = -1996.3334  # N4
$ws.Cells.Item(27, 8).Value = 37552  # H27
$ws.Cells.Item(27, 9).Value = 37549  # I27
$ws.Cells.Item(27, 10).Value = 37555  # J27
$ws.Cells.Item(27, 11).Value = 37549  # K27
$ws.Cells.Item(27, 12).Value = 37555  # L27
$ws.Cells.Item(27, 13).Value = -37357  # M27
$ws.Cells.Item(27, 14).Value = -37939  # N27
$ws.Cells.Item(70, 8).Value = 479000  # H70
$ws.Cells.Item(70, 10).Value = 479000  # J70
$ws.Cells.Item(70, 12).Value = 479000  # L70
$ws.Cells.Item(70, 14).Value = -479586  # N70
$ws.Cells.Item(73, 8).Value = 479000  # H73
$ws.Cells.Item(73, 10).Value = 479000  # J73
$ws.Cells.Item(73, 12).Value = 479000  # L73
$ws.Cells.Item(73, 14).Value = -481028  # N73
$ws.Cells.Item(134, 8).Value = 2103200.2  # H134
$ws.Cells.Item(134, 9).Value = 2383126.8  # I134
$ws.Cells.Item(134, 11).Value = 7149380.399999999  # K134
$ws.Cells.Item(134, 13).Value = -7146845.399999999  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2102.75  # H58
$ws.Cells.Item(58, 9).Value = 1757.4667  # I58
$ws.Cells.Item(58, 11).Value = 1757.4667  # K58
$ws.Cells.Item(58, 13).Value = -1554.4667  # M58
$ws.Cells.Item(122, 8).Value = 33336000  # H122
$ws.Cells.Item(122, 9).Value = 100000000  # I122
$ws.Cells.Item(122, 10).Value = 3999.5  # J122
$ws.Cells.Item(122, 11).Value = 300000000  # K122
$ws.Cells.Item(122, 12).Value = 11998.5  # L122
$ws.Cells.Item(122, 13).Value = -299997550  # M122
$ws.Cells.Item(122, 14).Value = -16898.5  # N122
$ws.Cells.Item(132, 8).Value = 2794.7097  # H132
$ws.Cells.Item(132, 9).Value = 2462  # I132
$ws.Cells.Item(132, 11).Value = 7386  # K132
$ws.Cells.Item(132, 13).Value = -4856  # M132
$ws.Cells.Item(134, 8).Value = 2718.678  # H134
$ws.Cells.Item(134, 9).Value = 1869.3469  # I134
$ws.Cells.Item(134, 10).Value = 6880.4  # J134
$ws.Cells.Item(134, 11).Value = 5608.0407  # K134
$ws.Cells.Item(134, 12).Value = 20641.2  # L134
$ws.Cells.Item(134, 13).Value = -3073.0407  # M134
$ws.Cells.Item(134, 14).Value = -25711.2  # N134
$ws.Cells.Item(136, 8).Value = 2102.75  # H136
$ws.Cells.Item(136, 9).Value = 1757.4667  # I136
$ws.Cells.Item(136, 11).Value = 5272.4001  # K136
$ws.Cells.Item(136, 13).Value = -2722.4001  # M136
$ws.Cells.Item(141, 8).Value = 585939.6  # H141
$ws.Cells.Item(141, 10).Value = 698641  # J141
$ws.Cells.Item(141, 12).Value = 698641  # L141
$ws.Cells.Item(141, 14).Value = -709001  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 14).ClearContents()  # N113
$ws.Cells.Item(113, 8).Value = 3149.5  # H113
$ws.Cells.Item(113, 9).Value = 3149.5  # I113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 11).Value = 9448.5  # K113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 13).Value = -7278.5  # M113

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 6859.4  # H122
$ws.Cells.Item(122, 9).Value = 6899  # I122
$ws.Cells.Item(122, 10).Value = 6849.5  # J122
$ws.Cells.Item(122, 11).Value = 20697  # K122
$ws.Cells.Item(122, 12).Value = 20548.5  # L122
$ws.Cells.Item(122, 13).Value = -18247  # M122
$ws.Cells.Item(122, 14).Value = -25448.5  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 900.5  # H16
$ws.Cells.Item(16, 9).Value = 828.0323  # I16
$ws.Cells.Item(16, 11).Value = 828.0323  # K16
$ws.Cells.Item(16, 13).Value = -658.0323  # M16
$ws.Cells.Item(55, 8).Value = 744.6111  # H55
$ws.Cells.Item(55, 9).Value = 845.4545000000001  # I55
$ws.Cells.Item(55, 11).Value = 845.4545000000001  # K55
$ws.Cells.Item(55, 13).Value = -672.4545000000001  # M55
$ws.Cells.Item(61, 8).Value = 3591.9546  # H61
$ws.Cells.Item(61, 9).Value = 2148.1177  # I61
$ws.Cells.Item(61, 11).Value = 2148.1177  # K61
$ws.Cells.Item(61, 13).Value = -1946.1177  # M61
$ws.Cells.Item(113, 8).Value = 3591.9546  # H113
$ws.Cells.Item(113, 9).Value = 2148.1177  # I113
$ws.Cells.Item(113, 11).Value = 2148.1177  # K113
$ws.Cells.Item(113, 13).Value = 21.88230000000021  # M113
$ws.Cells.Item(122, 8).Value = 6905952  # H122
$ws.Cells.Item(122, 9).Value = 17246880  # I122
$ws.Cells.Item(122, 11).Value = 51740640  # K122
$ws.Cells.Item(122, 13).Value = -51738190  # M122
$ws.Cells.Item(132, 8).Value = 530911  # H132
$ws.Cells.Item(132, 9).Value = 1114711.5  # I132
$ws.Cells.Item(132, 10).Value = 5490.5  # J132
$ws.Cells.Item(132, 11).Value = 3344134.5  # K132
$ws.Cells.Item(132, 12).Value = 16471.5  # L132
$ws.Cells.Item(132, 13).Value = -3341604.5  # M132
$ws.Cells.Item(132, 14).Value = -21531.5  # N132
$ws.Cells.Item(136, 8).Value = 5017.5557  # H136
$ws.Cells.Item(136, 9).Value = 4309  # I136
$ws.Cells.Item(136, 10).Value = 7497.5  # J136
$ws.Cells.Item(136, 11).Value = 12927  # K136
$ws.Cells.Item(136, 12).Value = 22492.5  # L136
$ws.Cells.Item(136, 13).Value = -10377  # M136
$ws.Cells.Item(136, 14).Value = -27592.5  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5840  # H62
$ws.Cells.Item(62, 9).Value = 3500  # I62
$ws.Cells.Item(62, 10).Value = 6425  # J62
$ws.Cells.Item(62, 11).Value = 3500  # K62
$ws.Cells.Item(62, 12).Value = 6425  # L62
$ws.Cells.Item(62, 13).Value = -2876  # M62
$ws.Cells.Item(62, 14).Value = -7673  # N62
$ws.Cells.Item(65, 8).Value = 5840  # H65
$ws.Cells.Item(65, 9).Value = 3500  # I65
$ws.Cells.Item(65, 10).Value = 6425  # J65
$ws.Cells.Item(65, 11).Value = 17500  # K65
$ws.Cells.Item(65, 12).Value = 32125  # L65
$ws.Cells.Item(65, 13).Value = -14380  # M65
$ws.Cells.Item(65, 14).Value = -38365  # N65
$ws.Cells.Item(126, 8).Value = 6543.3184  # H126
$ws.Cells.Item(126, 9).Value = 5167.8184  # I126
$ws.Cells.Item(126, 10).Value = 7918.8184  # J126
$ws.Cells.Item(126, 11).Value = 15503.4552  # K126
$ws.Cells.Item(126, 12).Value = 23756.4552  # L126
$ws.Cells.Item(126, 13).Value = -13033.4552  # M126
$ws.Cells.Item(126, 14).Value = -28696.4552  # N126
